$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C2:C6) from 2023-10-22 to 2023-10-25
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = (Get-Date -Year 2023 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0).Date
}
